$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 459.25
$ws.Range("I33").Value = 334.55554
$ws.Range("K33").Value = 334.55554
$ws.Range("M33").Value = -105.55554

$ws.Range("H64").Value = 3220
$ws.Range("I64").Value = 2990
$ws.Range("J64").Value = 3680
$ws.Range("K64").Value = 2990
$ws.Range("L64").Value = 3680
$ws.Range("M64").Value = -2742
$ws.Range("N64").Value = -4176

$ws.Range("H67").Value = 3220
$ws.Range("I67").Value = 2990
$ws.Range("J67").Value = 3680
$ws.Range("K67").Value = 2990
$ws.Range("L67").Value = 3680
$ws.Range("M67").Value = -2132
$ws.Range("N67").Value = -5396

$ws.Range("H138").Value = 2336.2163
$ws.Range("I138").Value = 1241
$ws.Range("J138").Value = 3002.8696
$ws.Range("K138").Value = 3723
$ws.Range("L138").Value = 9008.6088
$ws.Range("M138").Value = 1417
$ws.Range("N138").Value = -19288.6088

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("N37").ClearContents()

$ws.Range("H45").Value = 35492.07
$ws.Range("I45").Value = 56294.11
$ws.Range("K45").Value = 56294.11
$ws.Range("M45").Value = -55917.11

$ws.Range("H74").Value = 1076.1428
$ws.Range("I74").Value = 1038.64
$ws.Range("J74").Value = 1169.9
$ws.Range("K74").Value = 1038.64
$ws.Range("L74").Value = 1169.9
$ws.Range("M74").Value = -164.6400000000001
$ws.Range("N74").Value = -2917.9

$ws.Range("H77").Value = 1076.1428
$ws.Range("I77").Value = 1038.64
$ws.Range("J77").Value = 1169.9
$ws.Range("K77").Value = 5193.200000000001
$ws.Range("L77").Value = 5849.5
$ws.Range("M77").Value = -825.2000000000007
$ws.Range("N77").Value = -14585.5

$ws.Range("H122").Value = 1298.3182
$ws.Range("I122").Value = 1103.0625
$ws.Range("J122").Value = 1819
$ws.Range("K122").Value = 3309.1875
$ws.Range("L122").Value = 5457
$ws.Range("M122").Value = -859.1875
$ws.Range("N122").Value = -10357

$ws.Range("H132").Value = 1635828.5
$ws.Range("I132").Value = 1950.6666
$ws.Range("J132").Value = 4903584
$ws.Range("K132").Value = 5851.9998
$ws.Range("L132").Value = 14710752
$ws.Range("M132").Value = -3321.9998
$ws.Range("N132").Value = -14715812

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 50000856
$ws.Range("I107").Value = 55556340
$ws.Range("J107").Value = 1500
$ws.Range("K107").Value = 55556340
$ws.Range("L107").Value = 1500
$ws.Range("M107").Value = -55554420
$ws.Range("N107").Value = -5340

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 29412476
$ws.Range("I58").Value = 40000668
$ws.Range("J58").Value = 835
$ws.Range("K58").Value = 40000668
$ws.Range("L58").Value = 835
$ws.Range("M58").Value = -40000465
$ws.Range("N58").Value = -1241

$ws.Range("H62").Value = 3725.4443
$ws.Range("I62").Value = 2862.25
$ws.Range("J62").Value = 4416
$ws.Range("K62").Value = 2862.25
$ws.Range("L62").Value = 4416
$ws.Range("M62").Value = -2238.25
$ws.Range("N62").Value = -5664

$ws.Range("H65").Value = 3725.4443
$ws.Range("I65").Value = 2862.25
$ws.Range("J65").Value = 4416
$ws.Range("K65").Value = 14311.25
$ws.Range("L65").Value = 22080
$ws.Range("M65").Value = -11191.25
$ws.Range("N65").Value = -28320

$ws.Range("H135").Value = 49424.75
$ws.Range("J135").Value = 49424.75
$ws.Range("L135").Value = 49424.75
$ws.Range("N135").Value = -59564.75

$ws.Range("H136").Value = 29412476
$ws.Range("I136").Value = 40000668
$ws.Range("J136").Value = 835
$ws.Range("K136").Value = 120002004
$ws.Range("L136").Value = 2505
$ws.Range("M136").Value = -119999454
$ws.Range("N136").Value = -7605

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 20836968
$ws.Range("I5").Value = 31746394
$ws.Range("J5").Value = 9881.817999999999
$ws.Range("K5").Value = 95239182
$ws.Range("L5").Value = 29645.454
$ws.Range("M5").Value = -95239070
$ws.Range("N5").Value = -29869.454

$ws.Range("H122").Value = 12504627
$ws.Range("I122").Value = 48077156
$ws.Range("J122").Value = 6171.5674
$ws.Range("K122").Value = 432694404
$ws.Range("L122").Value = 55544.1066
$ws.Range("M122").Value = -432691954
$ws.Range("N122").Value = -60444.1066

$ws.Range("H131").Value = 768.27
$ws.Range("I131").Value = 455.9
$ws.Range("J131").Value = 802.9778
$ws.Range("K131").Value = 1367.7
$ws.Range("L131").Value = 2408.9334
$ws.Range("M131").Value = 3672.3
$ws.Range("N131").Value = -12488.9334

$ws.Range("H135").Value = 20836968
$ws.Range("I135").Value = 31746394
$ws.Range("J135").Value = 9881.817999999999
$ws.Range("K135").Value = 285717546
$ws.Range("L135").Value = 88936.36199999999
$ws.Range("M135").Value = -285715011
$ws.Range("N135").Value = -94006.36199999999

$ws.Range("H137").Value = 45455964
$ws.Range("I137").Value = 50000760
$ws.Range("J137").Value = 8000
$ws.Range("K137").Value = 150002280
$ws.Range("L137").Value = 24000
$ws.Range("M137").Value = -149997180
$ws.Range("N137").Value = -34200

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 22732680
$ws.Range("I122").Value = 35721956
$ws.Range("K122").Value = 107165868
$ws.Range("M122").Value = -107163418

$ws.Range("H132").Value = 4264.854
$ws.Range("I132").Value = 1937.8649
$ws.Range("J132").Value = 12092
$ws.Range("K132").Value = 5813.5947
$ws.Range("L132").Value = 36276
$ws.Range("M132").Value = -3283.5947
$ws.Range("N132").Value = -41336

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 9333.333000000001
$ws.Range("I122").Value = 9714.286
$ws.Range("K122").Value = 29142.858
$ws.Range("M122").Value = -26692.858

$ws.Range("H132").Value = 26413382
$ws.Range("I132").Value = 45705628
$ws.Range("J132").Value = 13465.895
$ws.Range("K132").Value = 137116884
$ws.Range("L132").Value = 40397.685
$ws.Range("M132").Value = -137114354
$ws.Range("N132").Value = -45457.685

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 1000
$ws.Range("I43").Value = 1000
$ws.Range("K43").Value = 1000
$ws.Range("M43").Value = -851

$ws.Range("H54").Value = 3000
$ws.Range("I54").Value = 3000
$ws.Range("K54").Value = 3000
$ws.Range("M54").Value = -2480

$ws.Range("H132").Value = 68523.234
$ws.Range("I132").Value = 104059.7
$ws.Range("J132").Value = 17756.857
$ws.Range("K132").Value = 312179.1
$ws.Range("L132").Value = 53270.571
$ws.Range("M132").Value = -309649.1
$ws.Range("N132").Value = -58330.571
